$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Title text: "PHIẾU DỊCH VỤ" -> "PHIẾU" / " THANH TOÁN" / " DỊCH VỤ"
# (three runs, identical bold/size formatting) with the "_GoBack"
# bookmark relocated in between " THANH TOÁN" and " DỊCH VỤ".
# Word only ever keeps a single "_GoBack" bookmark, so re-adding it at
# the new spot automatically removes it from its old location next to
# "Bác sĩ xác nhận" in the signature block.
# ------------------------------------------------------------------

$null = $d.Content.Find.Execute("PHIẾU DỊCH VỤ", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "PHIẾU THANH TOÁN DỊCH VỤ", 2)

$titleRange = $d.Paragraphs(1).Range
$titleStart = $titleRange.Start

# Force a run break right after "PHIẾU" (no bookmark should remain there,
# so drop a throw-away bookmark at that position and delete it again --
# deleting a bookmark removes its markup but leaves the run split intact).
$splitPoint1 = $d.Range($titleStart + 5, $titleStart + 5)
$d.Bookmarks.Add("ZZZTempSplit", $splitPoint1)
$d.Bookmarks("ZZZTempSplit").Delete()

# Force the second run break between " THANH TOÁN" and " DỊCH VỤ" by
# relocating the "_GoBack" bookmark to exactly that position.
$splitPoint2 = $d.Range($titleStart + 16, $titleStart + 16)
$d.Bookmarks.Add("_GoBack", $splitPoint2)
